# "create compare tab begining"
# 1) Rename the four existing sheets, dropping the "pointer_chase " prefix.
# 2) Round the "intensity" row (row 2) to whole numbers and rescale the
#    "energy" row (row 5) from micro-joules to milli-joules (divide by 1000)
#    on each of those sheets.
# 3) Add four new "compare" sheets at the end of the workbook (data_RAM
#    code_FLASH / data_RAM code_CCM / data_CCM code_FLASH / data_CCM
#    code_CCM) seeded with the same (now-rescaled) data as the four
#    existing sheets, one-to-one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: rename sheets 1-4
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "CF_DR"
$wb.Worksheets.Item(2).Name = "CC_DR"
$wb.Worksheets.Item(3).Name = "CF_DC"
$wb.Worksheets.Item(4).Name = "CC_DC"

# ---------------------------------------------------------------------
# Step 2: rewrite intensity (row 2) + energy (row 5) on each sheet
# ---------------------------------------------------------------------
$data = @{
    "CF_DR" = @{
        Intensity = @(16201, 22171, 27773)
        Energy    = @(10.826, 11.853, 17.32)
    }
    "CC_DR" = @{
        Intensity = @(13181, 23404, 32008)
        Energy    = @(8.816000000000001, 9.401999999999999, 10.026)
    }
    "CF_DC" = @{
        Intensity = @(14000, 23476, 31550)
        Energy    = @(9.368, 9.430999999999999, 9.875999999999999)
    }
    "CC_DC" = @{
        Intensity = @(11407, 22891, 33848)
        Energy    = @(7.621, 7.661, 7.575)
    }
}

$order = @("CF_DR", "CC_DR", "CF_DC", "CC_DC")

foreach ($name in $order) {
    $ws = $wb.Worksheets.Item($name)
    $vals = $data[$name]

    $ws.Range("B2").Value = $vals.Intensity[0]
    $ws.Range("C2").Value = $vals.Intensity[1]
    $ws.Range("D2").Value = $vals.Intensity[2]

    $ws.Range("B5").Value = $vals.Energy[0]
    $ws.Range("C5").Value = $vals.Energy[1]
    $ws.Range("D5").Value = $vals.Energy[2]
}

# ---------------------------------------------------------------------
# Step 3: add the four new "compare" sheets, copying each source sheet's
# full (now-updated) data across.
# ---------------------------------------------------------------------
$newSheets = @(
    @{ Name = "data_RAM code_FLASH"; Source = "CF_DR" },
    @{ Name = "data_RAM code_CCM";   Source = "CC_DR" },
    @{ Name = "data_CCM code_FLASH"; Source = "CF_DC" },
    @{ Name = "data_CCM code_CCM";   Source = "CC_DC" }
)

$cols = @("A", "B", "C", "D")
$rows = @(1, 2, 3, 4, 5)
# Row 1 (the "24"/"48"/"72" column headers) are stored as text in the
# source sheets even though they look numeric - force the same on copy.
$textRows = @(1)

foreach ($entry in $newSheets) {
    $src = $wb.Worksheets.Item($entry.Source)
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $new.Name = $entry.Name

    foreach ($r in $rows) {
        foreach ($c in $cols) {
            $addr = "$c$r"
            $cell = $src.Range($addr)
            $val = $cell.Value2
            # skip genuinely empty source cells (e.g. A1)
            if ($val -ne $null) {
                if ($textRows -contains $r) {
                    $new.Range($addr).NumberFormat = "@"
                    $new.Range($addr).Value = [string]$val
                } else {
                    $new.Range($addr).Value = $val
                }
            }
        }
    }
}

# Keep the originally-active first tab selected, matching the workbook's
# unchanged bookViews (activeTab="0").
$wb.Worksheets.Item(1).Activate()
